$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the text content (columns A-D) between row 42 and row 43,
# and move the Station18 (column J) value of 0 from row 43 to row 42.

$ws.Range("A42").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B42").Value = "unassigned"
$ws.Range("C42").Value = "unassigned"
$ws.Range("D42").Value = "unassigned"
$ws.Range("J42").Value = 0

$ws.Range("A43").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B43").Value = "Homo sapiens"
$ws.Range("C43").Value = "Human"
$ws.Range("D43").Value = "Human"
$ws.Range("J43").ClearContents()
